# Add a new "2023" column (T) to the 17.1.2 data table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data point: header year 2023 in T4, value 75.1 in T5, same styles as
# the existing year columns (S4/S5).
$ws.Range("T4").Value = 2023
$ws.Range("T4").Style = $ws.Range("S4").Style

$ws.Range("T5").Value = 75.1
$ws.Range("T5").Style = $ws.Range("S5").Style

# Slightly narrow the label columns (A:C) and set an explicit width for the
# newly-used data columns (D:T) to match the regenerated layout.
$ws.Range("A1:C1").EntireColumn.ColumnWidth = 35.5703125
$ws.Range("D1:T1").EntireColumn.ColumnWidth = 8.85546875

# Clear the lingering selection outside the used range (e.g. P8) by
# reselecting a cell inside the table.
$ws.Range("A1").Select()
